# Weekly Fruta/Hortalizas update: a new Granada price record (week of
# 2023-05-22) is inserted ahead of the existing rows, pushing the prior
# rows 36-45 down to 37-46.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 36 - this shifts the old rows
# 36..45 down to 37..46 (preserving their contents/formatting untouched).
$ws.Rows(36).Insert()

# Populate the newly inserted row 36 with the new weekly record.
$ws.Range("A36").Value = 9
$ws.Range("B36").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C36").Value = "Metropolitana"
$ws.Range("D36").Value = 45068
$ws.Range("E36").Value = 13
$ws.Range("F36").Value = "Fruta"
$ws.Range("G36").Value = 100104
$ws.Range("H36").Value = "Frutos de pepita"
$ws.Range("I36").Value = 100104001
$ws.Range("J36").Value = "Granada"
$ws.Range("K36").Value = "Wonderfull"
$ws.Range("L36").Value = "Primera"
$ws.Range("M36").Value = 350
$ws.Range("N36").Value = 10500
$ws.Range("O36").Value = 11000
$ws.Range("P36").Value = 10786
$ws.Range("Q36").Value = "$/caja 15 kilos granel"
$ws.Range("R36").Value = "Provincia de Curicó"
$ws.Range("S36").Value = 719
$ws.Range("T36").Value = 15
